$wb = $excel.ActiveWorkbook

# --- Main edit: add a new test-case row on the "Thread handling" sheet ---
# Before: column B (rows 3-7) had the style applied but no value, and row 8
# only had column A filled in. This commit marks all scenarios (including
# the newly added "Memory allocation" row) as supported ("Y"), and extends
# the table down to B8.
$ws = $wb.Worksheets.Item("Thread handling")
$ws.Activate()

# Fill existing empty cells B3:B7 with "Y"
$ws.Range("B3:B7").Value = "Y"

# Create the new B8 cell: copy formatting from B7 (style "s=5") then set value
$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "Y"

# Select B7:B8, matching the final selection left behind in the sheet view
$ws.Range("B7:B8").Select()

# --- "Package manager" sheet: selection grows from B8 to B7:B8 ---
$ws8 = $wb.Worksheets.Item("Package manager")
$ws8.Activate()
$ws8.Range("B7:B8").Select()

# Restore the originally active sheet/tab
$ws.Activate()
